$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.859.46"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "2.862.54"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.89%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.90%  "
$ws.Range("D9").Value = "2.869.08"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("D13").Value = "3.373.08"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "60.168.92"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.48%  "
$ws.Range("D17").Value = "2.887.61"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000136"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "345.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.436"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.11%  "
$ws.Range("E27").Value = "  -7.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.93%  "
$ws.Range("D30").Value = "0.0₃0833"
$ws.Range("E30").Value = "  -11.15%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.959"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.642"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").Value = "2.245.35"
$ws.Range("E41").Value = "  -7.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0567"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.58%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.16%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -12.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0893"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.51%  "
